{"js": "// Add the \"7 Mai\" entry (with its two sub-bullets) right after the\n// \"D\u00e9but de l'implementation SD\" bullet, at the very end of the document.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The anchor is the last paragraph in the document\n// (\"D\u00e9but de l'implementation SD\").\nconst anchor = paragraphs.items[paragraphs.items.length - 1];\n\n// Grab the list this paragraph belongs to, so the new paragraphs can be\n// attached to the very same numbering list (numId) instead of starting a\n// fresh one.\nconst anchorList = anchor.listOrNullObject;\nanchorList.load(\"id\");\nawait context.sync();\nconst listId = anchorList.id;\n\n// New top-level bullet: \"7 Mai :\" (same ilvl=0 as \"4 Mai :\").\nconst dateHeading = anchor.insertParagraph(\"7 Mai :\", \"After\");\ndateHeading.style = \"Paragraphedeliste\";\ndateHeading.attachToList(listId, 0);\n\n// Sub-bullet: \"Integration des capteurs, pas de tests\" (ilvl=1).\nconst sub1 = dateHeading.insertParagraph(\"Integration des capteurs, pas de tests\", \"After\");\nsub1.style = \"Paragraphedeliste\";\nsub1.attachToList(listId, 1);\n\n// Sub-bullet: \"Pas bc d'enerige aujourdhuis.\" (ilvl=1).\nconst sub2 = sub1.insertParagraph(\"Pas bc d\\u2019enerige aujourdhuis.\", \"After\");\nsub2.style = \"Paragraphedeliste\";\nsub2.attachToList(listId, 1);\n\nawait context.sync();\n", "ps1": "# Add the \"7 Mai\" entry (with its two sub-bullets) right after the\n# \"D\u00e9but de l'implementation SD\" bullet, at the very end of the document.\n$d = $word.ActiveDocument\n\n# Anchor on the last paragraph in the document\n# (\"D\u00e9but de l'implementation SD\").\n$anchor = $d.Paragraphs.Last\n\n# New top-level bullet: \"7 Mai :\" (same ilvl=0 / numId=1 list as \"4 Mai :\").\n$r1 = $anchor.Range\n$r1.Collapse(0)\n$r1.InsertParagraphAfter()\n$p1 = $d.Paragraphs.Item($d.Paragraphs.Count)\n$p1.Range.Text = \"7 Mai :\"\n$p1.Style = \"Paragraphedeliste\"\n$p1.Range.ListFormat.ListLevelNumber = 1\n\n# Sub-bullet: \"Integration des capteurs, pas de tests\" (ilvl=1).\n$r2 = $p1.Range\n$r2.Collapse(0)\n$r2.InsertParagraphAfter()\n$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)\n$p2.Range.Text = \"Integration des capteurs, pas de tests\"\n$p2.Style = \"Paragraphedeliste\"\n$p2.Range.ListFormat.ListLevelNumber = 2\n\n# Sub-bullet: \"Pas bc d'enerige aujourdhuis.\" (ilvl=1).\n$r3 = $p2.Range\n$r3.Collapse(0)\n$r3.InsertParagraphAfter()\n$p3 = $d.Paragraphs.Item($d.Paragraphs.Count)\n$p3.Range.Text = \"Pas bc d\" + [char]0x2019 + \"enerige aujourdhuis.\"\n$p3.Style = \"Paragraphedeliste\"\n$p3.Range.ListFormat.ListLevelNumber = 2\n"}
